# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# The source data pipeline swapped the stored order of 7 pairs of match
# rows in the "Portugal Primeira Liga" sheet. For every pair below, all of
# the match data in columns B (id) through AD (closing-line stats) needs
# to be exchanged between the two rows, while column A (the sequential
# match index) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(14, 15),
    @(48, 49),
    @(133, 134),
    @(167, 168),
    @(195, 196),
    @(282, 283),
    @(291, 292)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AD${r1}")
    $range2 = $ws.Range("B${r2}:AD${r2}")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
